# "restoration still buggy #1"
# - B3 value changes from "all" to "soe"
# - new row 5 added with a second restauration-related parameter
# - sheet view zoom bumped to 212% and selection left on B3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "soe"

$ws.Range("A5").Value = "second_amb_restauration"
$ws.Range("B5").Value = "yes"
$ws.Range("C5").Value = "yes or no. If yes, use restauration potential to restaurate ambitiously"

$excel.ActiveWindow.Zoom = 212
$ws.Range("B3").Select()
